{"js": "// Apply the four wording tweaks described by the diff:\n//   1. \"...important one which involves...\"      -> \"...important one that involves...\"\n//   2. \"...morphemes which are smallest...\"        -> \"...morphemes which are the smallest...\"\n//   3. \"...is segmentation which involves extraction of...\"\n//                                                    -> \"...is segmentation that involves the extraction of...\"\n//   4. \"...and other Africa languages...\"           -> \"...and other African languages...\"\n//\n// Each replacement is scoped to a unique surrounding phrase first (via\n// context.document.body.search) so that the correct occurrence of a common\n// word like \"which\" is targeted, then a narrower in-place search+replace\n// performs the actual text surgery without disturbing unrelated runs.\n\nconst body = context.document.body;\n\nasync function replaceWithin(scopeText, findText, replaceText) {\n  const scopeResults = body.search(scopeText, { matchCase: true });\n  scopeResults.load(\"items\");\n  await context.sync();\n\n  if (scopeResults.items.length === 0) {\n    throw new Error(`Scope text not found: ${scopeText}`);\n  }\n  const scope = scopeResults.items[0];\n\n  const hits = scope.search(findText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length === 0) {\n    throw new Error(`Find text not found within scope: ${findText}`);\n  }\n  hits.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"which\" -> \"that\" in \"...important one which involves splitting...\"\nawait replaceWithin(\n  \"segmentation is an important one which involves splitting structures such as words,\",\n  \"which\",\n  \"that\"\n);\n\n// 2) insert \"the \" before \"smallest meaning-bearing units of language\"\nawait replaceWithin(\n  \"s which are smallest meaning-bearing units of language\",\n  \"smallest meaning-bearing units of language\",\n  \"the smallest meaning-bearing units of language\"\n);\n\n// 3) \"which\" -> \"that\" and insert \"the \" before \"extraction\" in\n//    \"...is segmentation which involves extraction of underlying morphemes which are\"\nawait replaceWithin(\n  \"is segmentation which involves extraction of underlying morphemes which are\",\n  \"is segmentation which involves extraction\",\n  \"is segmentation that involves the extraction\"\n);\n\n// 4) \"Africa\" -> \"African\" in \"...and other Africa languages...\"\nawait replaceWithin(\n  \" and other Africa \",\n  \"Africa\",\n  \"African\"\n);\n", "ps1": "# Apply the four wording tweaks described by the diff:\n#   1. \"...important one which involves...\"        -> \"...important one that involves...\"\n#   2. \"...morphemes which are smallest...\"          -> \"...morphemes which are the smallest...\"\n#   3. \"...is segmentation which involves extraction of...\"\n#                                                      -> \"...is segmentation that involves the extraction of...\"\n#   4. \"...and other Africa languages...\"             -> \"...and other African languages...\"\n#\n# Each edit uses Find/Replace scoped with enough surrounding context to be\n# unique in the document (e.g. \"Africa\" alone would also match the\n# \"South Africa\" affiliation line in the title block, so the full\n# \" and other Africa \" phrase is used instead).\n\n$wdReplaceOne = 1\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text([string]$findText, [string]$replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, $wdReplaceOne)\n    if (-not $found) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# 1) \"which\" -> \"that\"\nReplace-Text \"which involves splitting\" \"that involves splitting\"\n\n# 2) insert \"the \" before \"smallest meaning-bearing units of language\"\nReplace-Text \"s which are smallest meaning-bearing units of language\" \"s which are the smallest meaning-bearing units of language\"\n\n# 3) \"which\" -> \"that\" and insert \"the \" before \"extraction\"\nReplace-Text \"is segmentation which involves extraction\" \"is segmentation that involves the extraction\"\n\n# 4) \"Africa\" -> \"African\"\nReplace-Text \" and other Africa \" \" and other African \"\n"}
